# ADD: sql query for days, dates and names
# Reassign the "Name" column (column B) values to the new order produced
# by the participant-rotation query, while Day (A) and Date (C) stay fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "4 Emrullah"
$ws.Range("B3").Value = "18 Kübra"
$ws.Range("B4").Value = "3 Eda"
$ws.Range("B5").Value = "6 Cahit"
$ws.Range("B6").Value = "8 Ertugrul"
$ws.Range("B7").Value = "10 Ümmü"
$ws.Range("B8").Value = "7 Nurullah"
$ws.Range("B9").Value = "29 Ahmet"
$ws.Range("B10").Value = "17 Sinan"
